$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-04-15 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-16 Tuesday", 2)

# Update the answer table cells in-place (row/column addressing avoids any
# text-collision issues and preserves each cell's existing paragraph/run
# formatting, since we are not creating new rows or paragraphs).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "47÷6=7, 5"
$t.Cell(1,2).Range.Text = "37÷8=4, 5"
$t.Cell(1,3).Range.Text = "98÷2=49, 0"
$t.Cell(1,4).Range.Text = "19÷6=3, 1"
$t.Cell(1,5).Range.Text = "12÷8=1, 4"

$t.Cell(5,1).Range.Text = "81÷9=9, 0"
$t.Cell(5,2).Range.Text = "33÷5=6, 3"
$t.Cell(5,3).Range.Text = "94÷4=23, 2"
$t.Cell(5,4).Range.Text = "98÷7=14, 0"
$t.Cell(5,5).Range.Text = "33÷3=11, 0"

$t.Cell(9,1).Range.Text = "85÷8=10, 5"
$t.Cell(9,2).Range.Text = "31÷4=7, 3"
$t.Cell(9,3).Range.Text = "49÷8=6, 1"
$t.Cell(9,4).Range.Text = "97÷5=19, 2"
$t.Cell(9,5).Range.Text = "40÷6=6, 4"

$t.Cell(13,1).Range.Text = "50÷4=12, 2"
$t.Cell(13,2).Range.Text = "22÷9=2, 4"
$t.Cell(13,3).Range.Text = "83÷7=11, 6"
$t.Cell(13,4).Range.Text = "27÷6=4, 3"
$t.Cell(13,5).Range.Text = "54÷4=13, 2"

$t.Cell(17,1).Range.Text = "76÷8=9, 4"
$t.Cell(17,2).Range.Text = "83÷7=11, 6"
$t.Cell(17,3).Range.Text = "27÷3=9, 0"
$t.Cell(17,4).Range.Text = "78÷4=19, 2"
$t.Cell(17,5).Range.Text = "95÷5=19, 0"
